# Auto-generated Excel COM-interop script to apply crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '22.080.63'
$ws.Range('E2').Value = '  +0.02%  '
$ws.Range('D3').Value = '1.559.43'
$ws.Range('E3').Value = '  +0.55%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9993'
$ws.Range('E4').Value = '  -0.24%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.9997'
$ws.Range('E5').Value = '  -0.24%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '291.80'
$ws.Range('E6').Value = '  +1.56%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3969'
$ws.Range('E7').Value = '  +4.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3245'
$ws.Range('E8').Value = '  -0.83%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '44.41'
$ws.Range('E9').Value = '  +2.21%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07290'
$ws.Range('E10').Value = '  -0.74%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.084'
$ws.Range('E11').Value = '  -3.78%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.9994'
$ws.Range('E12').Value = '  -0.22%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.727'
$ws.Range('E13').Value = '  -0.77%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '18.92'
$ws.Range('E14').Value = '  -5.65%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.668'
$ws.Range('E15').Value = '  -0.81%  '
$ws.Range('E16').Value = '  +5.16%  '
$ws.Range('D17').Value = '1.562.33'
$ws.Range('E17').Value = '  +0.39%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06593'
$ws.Range('E18').Value = '  -0.68%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '84.02'
$ws.Range('E19').Value = '  -1.83%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9988'
$ws.Range('E20').Value = '  -0.35%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.305'
$ws.Range('E21').Value = '  -0.54%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '15.64'
$ws.Range('E22').Value = '  -2.23%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.36'
$ws.Range('E23').Value = '  -2.32%  '
$ws.Range('D24').Value = '22.092.07'
$ws.Range('E24').Value = '  +0.06%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.367'
$ws.Range('E25').Value = '  +3.08%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.439'
$ws.Range('E26').Value = '  -2.33%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '148.67'
$ws.Range('E27').Value = '  -1.01%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.65'
$ws.Range('E28').Value = '  -2.38%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.868'
$ws.Range('E29').Value = '  -1.49%  '
$ws.Range('D30').Value = '1.732.82'
$ws.Range('E30').Value = '  -0.06%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '119.71'
$ws.Range('E31').Value = '  -1.43%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.9944'
$ws.Range('E32').Value = '  -7.15%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.912'
$ws.Range('E33').Value = '  +0.92%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08347'
$ws.Range('E34').Value = '  +1.83%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.177'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.606'
$ws.Range('E36').Value = '  -15.62%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02279'
$ws.Range('E37').Value = '  -1.59%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.156'
$ws.Range('E38').Value = '  -1.71%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06031'
$ws.Range('E39').Value = '  -3.89%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.211'
$ws.Range('E40').Value = '  -1.36%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.2057'
$ws.Range('E41').Value = '  -3.83%  '
$ws.Range('B42').Value = 'Aptos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '10.81'
$ws.Range('E42').Value = '  -1.49%  '
$ws.Range('B43').Value = 'Frax'
$ws.Range('C43').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.9990'
$ws.Range('E43').Value = '  -0.30%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.5849'
$ws.Range('E44').Value = '  -2.52%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.16'
$ws.Range('E45').Value = '  -3.47%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.758'
$ws.Range('E46').Value = '  +0.70%  '
$ws.Range('E47').Value = '  -3.21%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '118.86'
$ws.Range('E48').Value = '  -2.16%  '
$ws.Range('E49').Value = '  -2.90%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.144'
$ws.Range('E50').Value = '  -2.24%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06830'
$ws.Range('E51').Value = '  -2.41%  '
